$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.840.00'
$ws.Range("E2").Value = '  +3.00%  '
$ws.Range("D3").Value = '3.560.27'
$ws.Range("E3").Value = '  +2.03%  '
$ws.Range("E4").Value = '  +0.18%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.77'
$ws.Range("D5").Style = $ws.Range("B2").Style
$ws.Range("E5").Value = '  +2.27%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '186.33'
$ws.Range("D6").Style = $ws.Range("B2").Style
$ws.Range("E6").Value = '  +1.90%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.630'
$ws.Range("D7").Style = $ws.Range("B2").Style
$ws.Range("E7").Value = '  +2.76%  '
$ws.Range("D8").Value = '3.547.37'
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  +20.45%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.651'
$ws.Range("D11").Style = $ws.Range("B2").Style
$ws.Range("E11").Value = '  +2.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '54.45'
$ws.Range("D12").Style = $ws.Range("B2").Style
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000317'
$ws.Range("D13").Style = $ws.Range("B2").Style
$ws.Range("E13").Value = '  +6.15%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.46'
$ws.Range("D14").Style = $ws.Range("B2").Style
$ws.Range("E14").Value = '  +0.57%  '
$ws.Range("D15").Value = '4.127.62'
$ws.Range("E15").Value = '  +2.02%  '
$ws.Range("D16").Value = '70.890.98'
$ws.Range("E16").Value = '  +3.35%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.592.29'
$ws.Range("E17").Value = '  +3.30%  '
$ws.Range("B18").Value = 'Chainlink'
$ws.Range("C18").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.18'
$ws.Range("D18").Style = $ws.Range("B2").Style
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.47'
$ws.Range("D19").Style = $ws.Range("B2").Style
$ws.Range("E19").Value = '  +1.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '569.87'
$ws.Range("D20").Style = $ws.Range("B2").Style
$ws.Range("E20").Value = '  +5.70%  '
$ws.Range("E21").Value = '  +0.87%  '
$ws.Range("E22").Value = '  -0.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '17.63'
$ws.Range("D23").Style = $ws.Range("B2").Style
$ws.Range("E23").Value = '  -9.48%  '
$ws.Range("E24").Value = '  +4.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '4.90'
$ws.Range("D25").Style = $ws.Range("B2").Style
$ws.Range("E25").Value = '  -1.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '94.79'
$ws.Range("D26").Style = $ws.Range("B2").Style
$ws.Range("E26").Value = '  +0.71%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.24'
$ws.Range("D27").Style = $ws.Range("B2").Style
$ws.Range("E27").Value = '  +4.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.95'
$ws.Range("D28").Style = $ws.Range("B2").Style
$ws.Range("E28").Value = '  +1.72%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.15'
$ws.Range("D29").Style = $ws.Range("B2").Style
$ws.Range("E29").Value = '  +1.40%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '32.46'
$ws.Range("D30").Style = $ws.Range("B2").Style
$ws.Range("E30").Value = '  +3.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.21'
$ws.Range("D31").Style = $ws.Range("B2").Style
$ws.Range("E31").Value = '  +0.51%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.27'
$ws.Range("D32").Style = $ws.Range("B2").Style
$ws.Range("E32").Value = '  -2.07%  '
$ws.Range("E33").Value = '  +3.17%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '62.99'
$ws.Range("D34").Style = $ws.Range("B2").Style
$ws.Range("E34").Value = '  -1.77%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.41'
$ws.Range("D35").Style = $ws.Range("B2").Style
$ws.Range("E35").Value = '  +13.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '545.68'
$ws.Range("D36").Style = $ws.Range("B2").Style
$ws.Range("E36").Value = '  -4.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.413'
$ws.Range("D37").Style = $ws.Range("B2").Style
$ws.Range("E37").Value = '  +4.64%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.88'
$ws.Range("D38").Style = $ws.Range("B2").Style
$ws.Range("E38").Value = '  +0.32%  '
$ws.Range("E39").Value = '  +9.89%  '
$ws.Range("D40").Value = '0.0₃0801'
$ws.Range("E40").Value = '  +5.41%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.00'
$ws.Range("D41").Style = $ws.Range("B2").Style
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("D42").Value = '3.588.56'
$ws.Range("E42").Value = '  +12.26%  '
$ws.Range("B43").Value = 'Kaspa'
$ws.Range("C43").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.137'
$ws.Range("D43").Style = $ws.Range("B2").Style
$ws.Range("E43").Value = '  +3.38%  '
$ws.Range("B44").Value = 'Stacks'
$ws.Range("C44").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.45'
$ws.Range("D44").Style = $ws.Range("B2").Style
$ws.Range("E44").Value = '  +4.11%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0463'
$ws.Range("D45").Style = $ws.Range("B2").Style
$ws.Range("E45").Value = '  +6.44%  '
$ws.Range("E46").Value = '  +1.40%  '
$ws.Range("E47").Value = '  -1.44%  '
$ws.Range("E48").Value = '  +3.30%  '
$ws.Range("E49").Value = '  +2.95%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.50'
$ws.Range("D50").Style = $ws.Range("B2").Style
$ws.Range("E50").Value = '  +14.76%  '
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.000265'
$ws.Range("D51").Style = $ws.Range("B2").Style
$ws.Range("E51").Value = '  +17.23%  '
